$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains numeric-looking text values that must stay
# as literal text (e.g. "13.00", "0.00001110") instead of being auto-converted
# to numbers by Excel, which would strip formatting such as trailing zeros.
# Mark each target cell as Text before writing its new value.
$priceCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume %, and for two swapped
# rows, the coin name/link as well) cell by cell.
$ws.Range("D2").Value = "27.740.02"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "1.797.00"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "306.39"
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("D7").Value = "0.4958"
$ws.Range("E7").Value = "  -5.23%  "
$ws.Range("D8").Value = "0.3842"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "0.09264"
$ws.Range("E9").Value = "  +16.52%  "
$ws.Range("D10").Value = "1.091"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "40.50"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "6.264"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "20.38"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "1.798.19"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "7.160"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "0.00001110"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "92.11"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "0.06554"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").Value = "5.893"
$ws.Range("D23").Value = "27.771.85"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "10.95"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "2.225"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "156.36"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "20.48"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "2.001.29"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "2.395"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").Value = "126.15"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").Value = "1.051"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.509"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.605"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").Value = "0.06798"
$ws.Range("E35").Value = "  -6.29%  "
$ws.Range("D36").Value = "8.870"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").Value = "0.02290"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "0.2128"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "11.31"
$ws.Range("E39").Value = "  -7.34%  "
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("D41").Value = "0.6120"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "0.9997"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Value = "1.137"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "13.00"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "0.5841"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.274"
$ws.Range("E46").Value = "  -6.64%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.659"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").Value = "123.12"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").Value = "1.931"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "1.167"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").Value = "0.06698"
$ws.Range("E51").Value = "  -0.50%  "

# Restore default ("Normal") style on the price cells now that the text
# values are locked in, so no stray custom number format is left applied.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
